# Weekly update: insert a new price record (Acelga / Feria Lagunitas de
# Puerto Montt) above the existing row 300, shifting the rest of the
# table down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(300).Insert()

$ws.Cells.Item(300, 1).Value = 4
$ws.Cells.Item(300, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(300, 3).Value = 'Los Lagos'
$ws.Cells.Item(300, 4).Value = 45127
$ws.Cells.Item(300, 5).Value = 10
$ws.Cells.Item(300, 6).Value = 100112009
$ws.Cells.Item(300, 7).Value = 'Acelga'
$ws.Cells.Item(300, 8).Value = 'Sin especificar'
$ws.Cells.Item(300, 9).Value = 'Primera'
$ws.Cells.Item(300, 10).Value = 50
$ws.Cells.Item(300, 11).Value = 10000
$ws.Cells.Item(300, 12).Value = 10000
$ws.Cells.Item(300, 13).Value = 10000
$ws.Cells.Item(300, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(300, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(300, 16).Value = 833
$ws.Cells.Item(300, 17).Value = 12
$ws.Cells.Item(300, 18).Value = 'Hortaliza'

$ws.Cells.Item(300, 4).NumberFormat = $ws.Cells.Item(301, 4).NumberFormat
